$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S22").Copy($ws.Range("T22:Y22"))
$ws.Range("T22").Value = 5
$ws.Range("U22").Value = 19
$ws.Range("V22").Value = 3
$ws.Range("W22").Value = 9
$ws.Range("X22").Value = 2
$ws.Range("Y22").Value = 10

$ws.Range("S27").Copy($ws.Range("T27:Y27"))
$ws.Range("T27").Value = 17
$ws.Range("U27").Value = 7
$ws.Range("V27").Value = 8
$ws.Range("W27").Value = 4
$ws.Range("X27").Value = 9
$ws.Range("Y27").Value = 3

$ws.Range("S32").Copy($ws.Range("T32:Y32"))
$ws.Range("T32").Value = 7
$ws.Range("U32").Value = 17
$ws.Range("V32").Value = 4
$ws.Range("W32").Value = 8
$ws.Range("X32").Value = 3
$ws.Range("Y32").Value = 9

$ws.Range("S37").Copy($ws.Range("T37:Y37"))
$ws.Range("T37").Value = 13
$ws.Range("U37").Value = 11
$ws.Range("V37").Value = 7
$ws.Range("W37").Value = 5
$ws.Range("X37").Value = 6
$ws.Range("Y37").Value = 6

$ws.Range("S42").Copy($ws.Range("T42:Y42"))
$ws.Range("T42").Value = 2
$ws.Range("U42").Value = 22
$ws.Range("V42").Value = 1
$ws.Range("W42").Value = 11
$ws.Range("X42").Value = 1
$ws.Range("Y42").Value = 11

$ws.Range("S47").Copy($ws.Range("T47:Y47"))
$ws.Range("T47").Value = 15
$ws.Range("U47").Value = 9
$ws.Range("V47").Value = 7
$ws.Range("W47").Value = 5
$ws.Range("X47").Value = 8
$ws.Range("Y47").Value = 4

$ws.Range("S52").Copy($ws.Range("T52:Y52"))
$ws.Range("T52").Value = 3
$ws.Range("U52").Value = 21
$ws.Range("V52").Value = 2
$ws.Range("W52").Value = 10
$ws.Range("X52").Value = 1
$ws.Range("Y52").Value = 11

$ws.Range("S57").Copy($ws.Range("T57:Y57"))
$ws.Range("T57").Value = 22
$ws.Range("U57").Value = 2
$ws.Range("V57").Value = 12
$ws.Range("W57").Value = 0
$ws.Range("X57").Value = 10
$ws.Range("Y57").Value = 2

$ws.Range("S62").Copy($ws.Range("T62:Y62"))
$ws.Range("T62").Value = 8
$ws.Range("U62").Value = 16
$ws.Range("V62").Value = 4
$ws.Range("W62").Value = 8
$ws.Range("X62").Value = 4
$ws.Range("Y62").Value = 8

$ws.Range("S67").Copy($ws.Range("T67:Y67"))
$ws.Range("T67").Value = 3
$ws.Range("U67").Value = 21
$ws.Range("V67").Value = 1
$ws.Range("W67").Value = 11
$ws.Range("X67").Value = 2
$ws.Range("Y67").Value = 10

$ws.Range("S72").Copy($ws.Range("T72:Y72"))
$ws.Range("T72").Value = 23
$ws.Range("U72").Value = 1
$ws.Range("V72").Value = 11
$ws.Range("W72").Value = 1
$ws.Range("X72").Value = 12
$ws.Range("Y72").Value = 0

$ws.Range("S77").Copy($ws.Range("T77:Y77"))
$ws.Range("T77").Value = 20
$ws.Range("U77").Value = 4
$ws.Range("V77").Value = 10
$ws.Range("W77").Value = 2
$ws.Range("X77").Value = 10
$ws.Range("Y77").Value = 2

$ws.Range("S82").Copy($ws.Range("T82:Y82"))
$ws.Range("T82").Value = 19
$ws.Range("U82").Value = 5
$ws.Range("V82").Value = 10
$ws.Range("W82").Value = 2
$ws.Range("X82").Value = 9
$ws.Range("Y82").Value = 3

$ws.Range("S87").Copy($ws.Range("T87:Y87"))
$ws.Range("T87").Value = 23
$ws.Range("U87").Value = 1
$ws.Range("V87").Value = 11
$ws.Range("W87").Value = 1
$ws.Range("X87").Value = 12
$ws.Range("Y87").Value = 0

$ws.Range("S92").Copy($ws.Range("T92:Y92"))
$ws.Range("T92").Value = 14
$ws.Range("U92").Value = 10
$ws.Range("V92").Value = 7
$ws.Range("W92").Value = 5
$ws.Range("X92").Value = 7
$ws.Range("Y92").Value = 5

$ws.Range("S97").Copy($ws.Range("T97:Y97"))
$ws.Range("T97").Value = 20
$ws.Range("U97").Value = 4
$ws.Range("V97").Value = 8
$ws.Range("W97").Value = 4
$ws.Range("X97").Value = 12
$ws.Range("Y97").Value = 0

$ws.Range("S102").Copy($ws.Range("T102:Y102"))
$ws.Range("T102").Value = 5
$ws.Range("U102").Value = 19
$ws.Range("V102").Value = 3
$ws.Range("W102").Value = 9
$ws.Range("X102").Value = 2
$ws.Range("Y102").Value = 10

$ws.Range("S107").Copy($ws.Range("T107:Y107"))
$ws.Range("T107").Value = 22
$ws.Range("U107").Value = 2
$ws.Range("V107").Value = 12
$ws.Range("W107").Value = 0
$ws.Range("X107").Value = 10
$ws.Range("Y107").Value = 2

$ws.Range("S112").Copy($ws.Range("T112:Y112"))
$ws.Range("T112").Value = 2
$ws.Range("U112").Value = 22
$ws.Range("V112").Value = 1
$ws.Range("W112").Value = 11
$ws.Range("X112").Value = 1
$ws.Range("Y112").Value = 11

$ws.Range("S117").Copy($ws.Range("T117:Y117"))
$ws.Range("T117").Value = 17
$ws.Range("U117").Value = 7
$ws.Range("V117").Value = 8
$ws.Range("W117").Value = 4
$ws.Range("X117").Value = 9
$ws.Range("Y117").Value = 3

$ws.Range("S122").Copy($ws.Range("T122:Y122"))
$ws.Range("T122").Value = 9
$ws.Range("U122").Value = 15
$ws.Range("V122").Value = 5
$ws.Range("W122").Value = 7
$ws.Range("X122").Value = 4
$ws.Range("Y122").Value = 8

$ws.Range("S127").Copy($ws.Range("T127:Y127"))
$ws.Range("T127").Value = 13
$ws.Range("U127").Value = 11
$ws.Range("V127").Value = 8
$ws.Range("W127").Value = 4
$ws.Range("X127").Value = 5
$ws.Range("Y127").Value = 7

$ws.Range("S132").Copy($ws.Range("T132:Y132"))
$ws.Range("T132").Value = 10
$ws.Range("U132").Value = 14
$ws.Range("V132").Value = 5
$ws.Range("W132").Value = 7
$ws.Range("X132").Value = 5
$ws.Range("Y132").Value = 7

$ws.Range("S137").Copy($ws.Range("T137:Y137"))
$ws.Range("T137").Value = 5
$ws.Range("U137").Value = 19
$ws.Range("V137").Value = 3
$ws.Range("W137").Value = 9
$ws.Range("X137").Value = 2
$ws.Range("Y137").Value = 10

$ws.Range("S142").Copy($ws.Range("T142:Y142"))
$ws.Range("T142").Value = 15
$ws.Range("U142").Value = 9
$ws.Range("V142").Value = 10
$ws.Range("W142").Value = 2
$ws.Range("X142").Value = 5
$ws.Range("Y142").Value = 7

$ws.Range("S147").Copy($ws.Range("T147:Y147"))
$ws.Range("T147").Value = 2
$ws.Range("U147").Value = 22
$ws.Range("V147").Value = 0
$ws.Range("W147").Value = 12
$ws.Range("X147").Value = 2
$ws.Range("Y147").Value = 10

$ws.Range("S152").Copy($ws.Range("T152:Y152"))
$ws.Range("T152").Value = 17
$ws.Range("U152").Value = 7
$ws.Range("V152").Value = 6
$ws.Range("W152").Value = 6
$ws.Range("X152").Value = 11
$ws.Range("Y152").Value = 1

$ws.Range("S157").Copy($ws.Range("T157:Y157"))
$ws.Range("T157").Value = 7
$ws.Range("U157").Value = 17
$ws.Range("V157").Value = 3
$ws.Range("W157").Value = 9
$ws.Range("X157").Value = 4
$ws.Range("Y157").Value = 8

$ws.Range("S162").Copy($ws.Range("T162:Y162"))
$ws.Range("T162").Value = 12
$ws.Range("U162").Value = 12
$ws.Range("V162").Value = 8
$ws.Range("W162").Value = 4
$ws.Range("X162").Value = 4
$ws.Range("Y162").Value = 8

$ws.Range("S167").Copy($ws.Range("T167:Y167"))
$ws.Range("T167").Value = 23
$ws.Range("U167").Value = 1
$ws.Range("V167").Value = 12
$ws.Range("W167").Value = 0
$ws.Range("X167").Value = 11
$ws.Range("Y167").Value = 1

$ws.Range("S172").Copy($ws.Range("T172:Y172"))
$ws.Range("T172").Value = 6
$ws.Range("U172").Value = 18
$ws.Range("V172").Value = 3
$ws.Range("W172").Value = 9
$ws.Range("X172").Value = 3
$ws.Range("Y172").Value = 9

$ws.Range("S177").Copy($ws.Range("T177:Y177"))
$ws.Range("T177").Value = 8
$ws.Range("U177").Value = 16
$ws.Range("V177").Value = 2
$ws.Range("W177").Value = 10
$ws.Range("X177").Value = 6
$ws.Range("Y177").Value = 6

$ws.Range("S182").Copy($ws.Range("T182:Y182"))
$ws.Range("T182").Value = 13
$ws.Range("U182").Value = 11
$ws.Range("V182").Value = 6
$ws.Range("W182").Value = 6
$ws.Range("X182").Value = 7
$ws.Range("Y182").Value = 5

$ws.Range("S187").Copy($ws.Range("T187:Y187"))
$ws.Range("T187").Value = 1
$ws.Range("U187").Value = 23
$ws.Range("V187").Value = 1
$ws.Range("W187").Value = 11
$ws.Range("X187").Value = 0
$ws.Range("Y187").Value = 12

$ws.Range("S192").Copy($ws.Range("T192:Y192"))
$ws.Range("T192").Value = 3
$ws.Range("U192").Value = 21
$ws.Range("V192").Value = 2
$ws.Range("W192").Value = 10
$ws.Range("X192").Value = 1
$ws.Range("Y192").Value = 11

$ws.Range("S197").Copy($ws.Range("T197:Y197"))
$ws.Range("T197").Value = 22
$ws.Range("U197").Value = 2
$ws.Range("V197").Value = 10
$ws.Range("W197").Value = 2
$ws.Range("X197").Value = 12
$ws.Range("Y197").Value = 0

$ws.Range("S202").Copy($ws.Range("T202:Y202"))
$ws.Range("T202").Value = 11
$ws.Range("U202").Value = 13
$ws.Range("V202").Value = 5
$ws.Range("W202").Value = 7
$ws.Range("X202").Value = 6
$ws.Range("Y202").Value = 6

$ws.Range("S207").Copy($ws.Range("T207:Y207"))
$ws.Range("T207").Value = 20
$ws.Range("U207").Value = 4
$ws.Range("V207").Value = 11
$ws.Range("W207").Value = 1
$ws.Range("X207").Value = 9
$ws.Range("Y207").Value = 3

$ws.Range("S212").Copy($ws.Range("T212:Y212"))
$ws.Range("T212").Value = 22
$ws.Range("U212").Value = 2
$ws.Range("V212").Value = 11
$ws.Range("W212").Value = 1
$ws.Range("X212").Value = 11
$ws.Range("Y212").Value = 1

$ws.Range("S217").Copy($ws.Range("T217:Y217"))
$ws.Range("T217").Value = 14
$ws.Range("U217").Value = 10
$ws.Range("V217").Value = 8
$ws.Range("W217").Value = 4
$ws.Range("X217").Value = 6
$ws.Range("Y217").Value = 6

$ws.Range("S222").Copy($ws.Range("T222:Y222"))
$ws.Range("T222").Value = 6
$ws.Range("U222").Value = 18
$ws.Range("V222").Value = 4
$ws.Range("W222").Value = 8
$ws.Range("X222").Value = 2
$ws.Range("Y222").Value = 10

$ws.Range("S227").Copy($ws.Range("T227:Y227"))
$ws.Range("T227").Value = 11
$ws.Range("U227").Value = 13
$ws.Range("V227").Value = 6
$ws.Range("W227").Value = 6
$ws.Range("X227").Value = 5
$ws.Range("Y227").Value = 7

$ws.Range("S232").Copy($ws.Range("T232:Y232"))
$ws.Range("T232").Value = 10
$ws.Range("U232").Value = 14
$ws.Range("V232").Value = 8
$ws.Range("W232").Value = 4
$ws.Range("X232").Value = 2
$ws.Range("Y232").Value = 10
